# Stand Up Meeting - Week 7 update
# Fill in the Week 7 (Guillermo Toloza Guzman, rows 16-18) answers that were
# previously left blank: "¿Qué se hizo ayer?", "¿Qué se hará hoy?" and
# "¿Qué cosas se oponen?" for each weekday column (Lunes..Viernes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - "¿Qué se hizo ayer?"
$ws.Range("C16").Value = "Nada"
$ws.Range("D16").Value = "Se logro lo propuesto"
$ws.Range("E16").Value = "Diagrama de estructura"
$ws.Range("F16").Value = "Reunion semanal grupal"
$ws.Range("G16").Value = "Logré la meta"

# Row 17 - "¿Qué se hará hoy?"
$ws.Range("C17").Value = "Coherencia diagramas"
$ws.Range("D17").Value = "Verificar cambios realizados y anexar informacion necesaria para la coherencia de diagramas"
$ws.Range("E17").Value = "Reunion semanal grupal"
$ws.Range("F17").Value = "Buscar sobre relacion muchos a muchos en el modelo de dominio"
$ws.Range("G17").Value = "Nada"

# Row 18 - "¿Qué cosas se oponen?"
$ws.Range("C18").Value = "Nada"
$ws.Range("D18").Value = "Nada"
$ws.Range("E18").Value = "Nada"
$ws.Range("F18").Value = "Nada"
$ws.Range("G18").Value = "Nada"

# Leave the selection where the author ended up while filling the sheet.
$ws.Range("G22").Select()
